$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add P1 = 14, Q1 = 15, copying the header style from O1 ---
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# --- Data rows 2-25: swap I<->K and M<->O, then add P=2, Q=2 ---
for ($r = 2; $r -le 25; $r++) {
    $iVal = $ws.Cells.Item($r, 9).Value2
    $kVal = $ws.Cells.Item($r, 11).Value2
    $ws.Cells.Item($r, 9).Value = $kVal
    $ws.Cells.Item($r, 11).Value = $iVal

    $mVal = $ws.Cells.Item($r, 13).Value2
    $oVal = $ws.Cells.Item($r, 15).Value2
    $ws.Cells.Item($r, 13).Value = $oVal
    $ws.Cells.Item($r, 15).Value = $mVal

    $ws.Cells.Item($r, 16).Value = 2
    $ws.Cells.Item($r, 17).Value = 2
}
